$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1122.3636
$ws.Range("I39").Value = 1292.6666
$ws.Range("K39").Value = 3877.9998
$ws.Range("M39").Value = -3581.9998
$ws.Range("H40").Value = 3241.2856
$ws.Range("I40").Value = 9800
$ws.Range("J40").Value = 2148.1667
$ws.Range("K40").Value = 9800
$ws.Range("L40").Value = 2148.1667
$ws.Range("M40").Value = -9625
$ws.Range("N40").Value = -2498.1667
$ws.Range("H98").Value = 3682.0588
$ws.Range("I98").Value = 2340.3125
$ws.Range("K98").Value = 2340.3125
$ws.Range("M98").Value = -842.3125
$ws.Range("H113").Value = 2938.4285
$ws.Range("J113").Value = 3350.5715
$ws.Range("L113").Value = 3350.5715
$ws.Range("N113").Value = -9858.5715
$ws.Range("H116").Value = 3079.4
$ws.Range("I116").Value = 2737.9565
$ws.Range("K116").Value = 2737.9565
$ws.Range("M116").Value = 704.0435000000002
$ws.Range("H122").Value = 3682.0588
$ws.Range("I122").Value = 2340.3125
$ws.Range("K122").Value = 7020.9375
$ws.Range("M122").Value = -4570.9375
$ws.Range("H136").Value = 15000
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1734.2727
$ws.Range("I61").Value = 1266.25
$ws.Range("J61").Value = 2001.7142
$ws.Range("K61").Value = 1266.25
$ws.Range("L61").Value = 2001.7142
$ws.Range("M61").Value = -1054.25
$ws.Range("N61").Value = -2425.7142
$ws.Range("H97").Value = 474.92856
$ws.Range("I97").Value = 445.75
$ws.Range("K97").Value = 445.75
$ws.Range("M97").Value = 50.25
$ws.Range("H102").Value = 9806807
$ws.Range("I102").Value = 11907629
$ws.Range("K102").Value = 11907629
$ws.Range("M102").Value = -11906007
$ws.Range("H110").Value = 1140.2858
$ws.Range("I110").Value = 908.86664
$ws.Range("J110").Value = 1718.8334
$ws.Range("K110").Value = 908.86664
$ws.Range("L110").Value = 1718.8334
$ws.Range("M110").Value = 1136.13336
$ws.Range("N110").Value = -5808.8334
$ws.Range("H122").Value = 1193.3
$ws.Range("I122").Value = 994.125
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 2982.375
$ws.Range("L122").Value = 5970
$ws.Range("M122").Value = -532.375
$ws.Range("N122").Value = -10870
$ws.Range("H132").Value = 2048.7715
$ws.Range("I132").Value = 1750.5555
$ws.Range("J132").Value = 3055.25
$ws.Range("K132").Value = 5251.666499999999
$ws.Range("L132").Value = 9165.75
$ws.Range("M132").Value = -2721.666499999999
$ws.Range("N132").Value = -14225.75
$ws.Range("H136").Value = 1734.2727
$ws.Range("I136").Value = 1266.25
$ws.Range("J136").Value = 2001.7142
$ws.Range("K136").Value = 3798.75
$ws.Range("L136").Value = 6005.142599999999
$ws.Range("M136").Value = -1248.75
$ws.Range("N136").Value = -11105.1426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 634.7143
$ws.Range("I80").Value = 409.2
$ws.Range("K80").Value = 409.2
$ws.Range("M80").Value = 588.8
$ws.Range("H83").Value = 634.7143
$ws.Range("I83").Value = 409.2
$ws.Range("K83").Value = 2046
$ws.Range("M83").Value = 2946
$ws.Range("H86").Value = 3834.5386
$ws.Range("I86").Value = 4320
$ws.Range("J86").Value = 3172.5454
$ws.Range("K86").Value = 4320
$ws.Range("L86").Value = 3172.5454
$ws.Range("M86").Value = -3197
$ws.Range("N86").Value = -5418.5454
$ws.Range("H89").Value = 3834.5386
$ws.Range("I89").Value = 4320
$ws.Range("J89").Value = 3172.5454
$ws.Range("K89").Value = 21600
$ws.Range("L89").Value = 15862.727
$ws.Range("M89").Value = -15984
$ws.Range("N89").Value = -27094.727
$ws.Range("H99").Value = 38462900
$ws.Range("I99").Value = 55556780
$ws.Range("K99").Value = 55556780
$ws.Range("M99").Value = -55555282
$ws.Range("H105").Value = 111114230
$ws.Range("I105").Value = 125003260
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 125003260
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -125001513
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 5088.759
$ws.Range("I134").Value = 982.96
$ws.Range("J134").Value = 30750
$ws.Range("K134").Value = 2948.88
$ws.Range("L134").Value = 92250
$ws.Range("M134").Value = -413.8800000000001
$ws.Range("N134").Value = -97320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1195.7216
$ws.Range("I31").Value = 1036.6571
$ws.Range("J31").Value = 2432.889
$ws.Range("K31").Value = 1036.6571
$ws.Range("L31").Value = 2432.889
$ws.Range("M31").Value = -741.6570999999999
$ws.Range("N31").Value = -3022.889
$ws.Range("H34").Value = 1195.7216
$ws.Range("I34").Value = 1036.6571
$ws.Range("J34").Value = 2432.889
$ws.Range("K34").Value = 1036.6571
$ws.Range("L34").Value = 2432.889
$ws.Range("M34").Value = -834.6570999999999
$ws.Range("N34").Value = -2836.889
$ws.Range("H58").Value = 1148
$ws.Range("I58").Value = 968.4286
$ws.Range("J58").Value = 1567
$ws.Range("K58").Value = 968.4286
$ws.Range("L58").Value = 1567
$ws.Range("M58").Value = -765.4286
$ws.Range("N58").Value = -1973
$ws.Range("H132").Value = 3016.4375
$ws.Range("I132").Value = 2568.875
$ws.Range("J132").Value = 3464
$ws.Range("K132").Value = 7706.625
$ws.Range("L132").Value = 10392
$ws.Range("M132").Value = -5176.625
$ws.Range("N132").Value = -15452
$ws.Range("H136").Value = 1148
$ws.Range("I136").Value = 968.4286
$ws.Range("J136").Value = 1567
$ws.Range("K136").Value = 2905.2858
$ws.Range("L136").Value = 4701
$ws.Range("M136").Value = -355.2857999999997
$ws.Range("N136").Value = -9801

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 357.81818
$ws.Range("J23").Value = 239.46666
$ws.Range("L23").Value = 718.3999799999999
$ws.Range("N23").Value = -1188.39998
$ws.Range("H34").Value = 1920.1
$ws.Range("I34").Value = 974
$ws.Range("J34").Value = 2429.5386
$ws.Range("K34").Value = 2922
$ws.Range("L34").Value = 7288.6158
$ws.Range("M34").Value = -2838
$ws.Range("N34").Value = -7456.6158
$ws.Range("H39").Value = 1938.3636
$ws.Range("J39").Value = 1644.421
$ws.Range("L39").Value = 4933.263
$ws.Range("N39").Value = -5521.263
$ws.Range("H55").Value = 2500
$ws.Range("J55").Value = 3250
$ws.Range("L55").Value = 9750
$ws.Range("N55").Value = -10104
$ws.Range("H113").Value = 609
$ws.Range("J113").Value = 609
$ws.Range("L113").Value = 1827
$ws.Range("N113").Value = -6167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2747.2942
$ws.Range("I102").Value = 2905.875
$ws.Range("K102").Value = 2905.875
$ws.Range("M102").Value = -1283.875
$ws.Range("H122").Value = 1437.2142
$ws.Range("I122").Value = 1525.1818
$ws.Range("J122").Value = 1114.6666
$ws.Range("K122").Value = 4575.5454
$ws.Range("L122").Value = 3343.9998
$ws.Range("M122").Value = -2125.5454
$ws.Range("N122").Value = -8243.9998
$ws.Range("H132").Value = 1826.069
$ws.Range("I132").Value = 1497.9
$ws.Range("K132").Value = 4493.700000000001
$ws.Range("M132").Value = -1963.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 998.5333000000001
$ws.Range("I16").Value = 1007.75
$ws.Range("J16").Value = 961.6667
$ws.Range("K16").Value = 1007.75
$ws.Range("L16").Value = 961.6667
$ws.Range("M16").Value = -837.75
$ws.Range("N16").Value = -1301.6667
$ws.Range("H93").Value = 979.7143
$ws.Range("I93").Value = 623.125
$ws.Range("J93").Value = 2120.8
$ws.Range("K93").Value = 623.125
$ws.Range("L93").Value = 2120.8
$ws.Range("M93").Value = 624.875
$ws.Range("N93").Value = -4616.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 339.125
$ws.Range("I107").Value = 258.83334
$ws.Range("J107").Value = 580
$ws.Range("K107").Value = 776.5000200000001
$ws.Range("L107").Value = 1740
$ws.Range("M107").Value = 1143.49998
$ws.Range("N107").Value = -5580
$ws.Range("H132").Value = 1383.7046
$ws.Range("I132").Value = 1103.75
$ws.Range("J132").Value = 2130.25
$ws.Range("K132").Value = 3311.25
$ws.Range("L132").Value = 6390.75
$ws.Range("M132").Value = -781.25
$ws.Range("N132").Value = -11450.75
$ws.Range("H136").Value = 445.19232
$ws.Range("I136").Value = 348.68182
$ws.Range("J136").Value = 976
$ws.Range("K136").Value = 1046.04546
$ws.Range("L136").Value = 2928
$ws.Range("M136").Value = 1503.95454
$ws.Range("N136").Value = -8028
